# BOM Expansion with three pos: temperature sensor, heating cable and relay
#
# Adds three new rows (Pos 10/11/12) to the "Electric Hardware" sheet's
# table (Temperatur Sensor, Relay, Heizkabel — qty 1 each), grows the
# table/dimension to cover them, and makes "Electric Hardware" the active
# sheet/tab with the selection parked on F18 (matching where the next free
# row after the expanded table would be).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Electric Hardware")

# --- new row 15: Pos 10, Temperatur Sensor, Anzahl 1 -----------------------
$ws.Cells.Item(15, 3).Value = 10
$ws.Cells.Item(15, 3).HorizontalAlignment = -4131   # xlLeft, matches "Pos" column style
$ws.Cells.Item(15, 4).Value = "Temperatur Sensor"
$ws.Cells.Item(15, 5).Value = 1

# --- new row 16: Pos 11, Relay, Anzahl 1 ------------------------------------
$ws.Cells.Item(16, 3).Value = 11
$ws.Cells.Item(16, 3).HorizontalAlignment = -4131
$ws.Cells.Item(16, 4).Value = "Relay"
$ws.Cells.Item(16, 5).Value = 1

# --- new row 17: Pos 12, Heizkabel, Anzahl 1 --------------------------------
$ws.Cells.Item(17, 3).Value = 12
$ws.Cells.Item(17, 3).HorizontalAlignment = -4131
$ws.Cells.Item(17, 4).Value = "Heizkabel"
$ws.Cells.Item(17, 5).Value = 1

# Grow the worksheet table ("Tabelle2") so it covers the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("C5:G17"))

# Electric Hardware becomes the active/selected sheet & tab.
$ws.Activate() | Out-Null
$ws.Range("F18").Select() | Out-Null
